$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.512.41'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.940.64'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.23%  '
$ws.Range("E9").Value = '  -3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0852'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("D12").Value = '2.227.18'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("E13").Value = '  -2.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.811'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.40%  '
$ws.Range("D17").Value = '1.934.96'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = '36.438.63'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("D20").Value = '0.0₃0864'
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.95%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.134'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.71%  '
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.09'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0613'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.98%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0982'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0208'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = '1.340.36'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").Value = '2.118.51'
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.42%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
